$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
Write-Output ($s | Get-Member -Name Design | Out-String)
